# Insert 3 new accelerometer readings at the top of the data (rows 2-4),
# pushing the existing data down. Because the sheet only grows by a net
# of 1 row (dimension goes from C20 to C21), the two rows that get
# pushed past the new bottom of the range are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 2, shifting existing data down.
$ws.Range("A2:C4").EntireRow.Insert()

# New rows of accelerometer data (x, y, z)
$newRows = @(
    @(-2.669419974088668, 9.347340643405914, -0.05590170621871929),
    @(-2.789929866790771, 9.389312267303467, -0.0143058076500895),
    @(-3.052737355232238, 9.127021908760071, -0.5663906224071975)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}

# The data that was previously on rows 19 and 20 has now been shifted
# down to rows 22 and 23, beyond the new dimension of A1:C21. Remove
# those trailing rows so the sheet ends at row 21.
$ws.Range("A22:C23").EntireRow.Delete()
